# Rename the "CLASSIFICATION" header (column C) to "INFORMATIONDOMAIN", and
# swap the header-row highlight formatting between columns B and C so the
# distinct formatting that used to mark the CLASSIFICATION cell now marks
# the renamed INFORMATIONDOMAIN cell (and the ORGANIZATION column reverts
# to the plain formatting that CLASSIFICATION used to have).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-CellFormat($range) {
    $font = $range.Font
    return @{
        FontName      = $font.Name
        FontSize      = $font.Size
        FontBold      = $font.Bold
        FontItalic    = $font.Italic
        FontColor     = $font.Color
        FontUnderline = $font.Underline
        HorizontalAlign = $range.HorizontalAlignment
        VerticalAlign    = $range.VerticalAlignment
        WrapText         = $range.WrapText
    }
}

function Set-CellFormat($range, $fmt) {
    $range.Font.Name = $fmt.FontName
    $range.Font.Size = $fmt.FontSize
    $range.Font.Bold = $fmt.FontBold
    $range.Font.Italic = $fmt.FontItalic
    $range.Font.Color = $fmt.FontColor
    $range.Font.Underline = $fmt.FontUnderline
    $range.HorizontalAlignment = $fmt.HorizontalAlign
    $range.VerticalAlignment = $fmt.VerticalAlign
    $range.WrapText = $fmt.WrapText
}

$rangeB1 = $ws.Range("B1")
$rangeC1 = $ws.Range("C1")
$rangeB2 = $ws.Range("B2")

# Capture the current (pre-edit) formatting of B1 and C1 before changing
# anything, so the swap below is based on the original state.
$formatB1 = Get-CellFormat $rangeB1
$formatC1 = Get-CellFormat $rangeC1

# Rename the shared text used by the header cell in column C.
$rangeC1.Value = "INFORMATIONDOMAIN"

# Swap the header-row formatting between columns B and C.
Set-CellFormat $rangeB1 $formatC1
Set-CellFormat $rangeC1 $formatB1

# Row 2's column-B cell also takes on the formatting that used to belong
# to column C.
Set-CellFormat $rangeB2 $formatC1
